# edit.ps1 -- apply the project_6.docx content overhaul (AI presentation-tool
# copy replacing the old India EV-market copy) plus the removal of the old
# "Page 2 - Section 3" block and the trailing paragraph of old Section/Page 2-1.
#
# Strategy: walk the (stable) paragraph collection by 1-based index and
# overwrite Range.Text in place so paragraph/run formatting (pStyle, sz, etc.)
# survives untouched; then delete the paragraphs that the edit drops entirely.
# Range.Text assignment is used instead of Find.Execute's replacement so the
# host's AutoCorrect-style smart-quote substitution never touches the new copy
# (it would silently turn straight "'" into a curly ’ otherwise).

$d = $word.ActiveDocument

function Set-ParaText($index, $expectedOld, $newText) {
    $p = $d.Paragraphs.Item($index)
    $cur = $p.Range.Text
    if ($cur -notmatch [regex]::Escape($expectedOld)) {
        Write-Output ("WARNING paragraph " + $index + " text was not as expected before replace: " + $cur)
    }
    $p.Range.Text = $newText
}

# --- Text replacements (title, headings, body paragraphs) ---
Set-ParaText 1 "Electric Vehicle Market in India 2025" "buienss ui"
Set-ParaText 2 "Page 1 – Section 1" "Section 1"
Set-ParaText 3 "The Indian electric vehicle (EV) market is poised for significant expansion by 2025, driven by increasing environmental awareness, government support, and declining battery costs. While still a nascent market compared to global leaders, India presents a substantial growth opportunity. Current estimates project a market size exceeding `$7 billion by 2025, representing a CAGR of over 40% from 2020 levels. This growth is expected across various vehicle segments, including two-wheelers, three-wheelers, passenger vehicles, and commercial vehicles." "Our AI-powered presentation tool is designed to revolutionize how presentations are created and delivered. By leveraging advanced artificial intelligence, the tool automates key aspects of the presentation process, from content generation to design optimization, making it accessible to users of all skill levels."
Set-ParaText 4 "The Government of India's initiatives, such as the Faster Adoption and Manufacturing of Electric Vehicles (FAME) scheme, play a crucial role in incentivizing EV adoption through subsidies and tax benefits. These policies, coupled with state-level initiatives promoting EV manufacturing and infrastructure development, are creating a favorable ecosystem for EV market growth. Furthermore, increasing fuel prices and growing concerns about air pollution in major cities are pushing consumers towards electric alternatives." "The core technology utilizes a combination of natural language processing (NLP), machine learning (ML), and computer vision to understand user inputs, generate relevant content suggestions, and create visually appealing slides. The tool is intended to democratize access to professional-quality presentations, regardless of the user's design or writing abilities."
Set-ParaText 5 "Despite the positive outlook, challenges remain. These include the high upfront cost of EVs, limited charging infrastructure, and range anxiety among consumers. Addressing these challenges through innovative business models, technological advancements, and strategic partnerships will be critical to realizing the full potential of the Indian EV market." "Key features include automated content generation, intelligent slide design, real-time feedback and suggestions, and seamless integration with popular presentation platforms. Users can input a topic, keyword, or existing document, and the tool will generate a complete presentation outline, including slide content, images, and design elements. The AI continuously learns from user feedback and data to improve its performance and provide increasingly relevant and effective presentations."
Set-ParaText 6 "Page 1 – Section 2" "Section 2"
Set-ParaText 7 "Several key factors are propelling the growth of the EV market in India. Firstly, government policies, including subsidies under the FAME scheme and tax exemptions, significantly reduce the total cost of ownership for EVs. Secondly, increasing environmental concerns and stricter emission norms are driving demand for cleaner transportation alternatives." "Our primary target audience comprises both students and business professionals. For students, the tool offers a significant advantage in creating compelling presentations for academic assignments, group projects, and research presentations, reducing the time spent on design and enabling them to focus on the content itself."
Set-ParaText 8 "Thirdly, falling battery prices are making EVs more affordable and competitive with traditional vehicles. Battery costs, which constitute a significant portion of the EV's price, have declined substantially in recent years and are expected to continue decreasing. Fourthly, rising fuel prices are making EVs a more attractive option for consumers seeking to reduce their transportation expenses." "For business professionals, the tool streamlines the presentation creation process for sales pitches, investor updates, internal training, and marketing materials. It enables teams to rapidly generate high-quality presentations, ensuring consistent branding and messaging across the organization."
Set-ParaText 9 "Finally, growing awareness among consumers about the benefits of EVs, coupled with increasing availability of EV models and charging infrastructure, is contributing to the market's expansion." "Specifically within the business segment, we are targeting small-to-medium sized businesses (SMBs) that often lack dedicated design resources or the time to create professional-looking presentations. Furthermore, sales and marketing teams within larger corporations are a key target due to their frequent need for impactful presentation materials."
Set-ParaText 10 "Page 1 – Section 3" "Section 3"
Set-ParaText 11 "The Indian EV market can be segmented based on vehicle type, application, and geography. By vehicle type, the market is divided into two-wheelers, three-wheelers, passenger vehicles, and commercial vehicles. Two-wheelers and three-wheelers currently dominate the market due to their affordability and suitability for urban transportation." "Key features of the tool include AI-powered content generation, which automatically creates slide content based on user-provided keywords or topics. The intelligent slide design feature optimizes layouts, selects appropriate visuals, and ensures consistent branding."
Set-ParaText 12 "By application, the market is segmented into personal and commercial use. The commercial segment is expected to witness significant growth, driven by the increasing adoption of electric buses and delivery vehicles. Geographically, the market is concentrated in metropolitan areas and Tier I cities, but is gradually expanding to Tier II and Tier III cities with increasing infrastructure development." "Real-time feedback and suggestions offer users immediate guidance on improving their presentations, covering aspects like content clarity, visual appeal, and delivery techniques. The platform also offers seamless integration with popular presentation software like PowerPoint and Google Slides, allowing users to easily export and share their creations."
Set-ParaText 13 "Further segmentation can be based on battery type (Lithium-ion, Lead-acid, etc.) and motor type (AC induction, Permanent magnet synchronous motor, etc.). Understanding these segmentations is crucial for companies to tailor their products and services to specific customer needs and market opportunities." "Furthermore, the tool includes a library of pre-designed templates and assets, providing users with a starting point for their presentations and allowing them to quickly customize existing designs. These templates cover a wide range of presentation styles and industries, catering to diverse user needs."
Set-ParaText 14 "Page 2 – Section 1" "Section 4"
Set-ParaText 15 "Technological advancements in battery technology, charging infrastructure, and vehicle components are shaping the EV landscape in India. Solid-state batteries, with their higher energy density and improved safety, are expected to become more prevalent in the future." "For students, our tool enhances learning outcomes by facilitating the creation of engaging and informative presentations. It reduces presentation anxiety by automating the design process, allowing students to focus on research and content mastery. The automated design assistance promotes better visual communication skills, beneficial for academic success and future career prospects."
Set-ParaText 16 "Wireless charging technology is also gaining traction, offering a convenient and user-friendly charging experience. Furthermore, advancements in motor technology are improving the efficiency and performance of EVs. The development of robust and reliable charging infrastructure, including both public and private charging stations, is critical for supporting the widespread adoption of EVs." "Businesses benefit from increased efficiency and productivity. Rapid presentation creation saves valuable time and resources. Consistent branding and messaging across all presentations enhances brand image and strengthens communication with stakeholders. Data-driven insights from user engagement metrics allow businesses to optimize their presentation strategy and improve audience impact."
Set-ParaText 18 "Page 2 – Section 2" "Section 5"
Set-ParaText 19 "Several factors influence purchasing behavior in the Indian EV market. Price sensitivity is a major consideration for most consumers, with the upfront cost of EVs being a significant barrier. Government subsidies and incentives play a crucial role in mitigating this barrier." "Our future roadmap includes enhancements to the AI engine, incorporating more sophisticated algorithms for content generation and design optimization. We plan to introduce personalized learning paths for users, tailored to their individual skill levels and presentation needs."
Set-ParaText 20 "Range anxiety is another important factor, with consumers concerned about the limited range of EVs and the availability of charging infrastructure. Trust in brands and the perceived reliability of EVs also influence purchasing decisions. Word-of-mouth referrals and online reviews play a significant role." "Integration with additional platforms and data sources is also a priority, expanding the tool's compatibility and enabling users to incorporate data-driven insights into their presentations more easily. We will also be exploring the addition of features such as automated translation and real-time collaboration tools to further enhance the user experience."
Set-ParaText 21 "Consumers are also increasingly concerned about the environmental impact of their transportation choices and are willing to pay a premium for eco-friendly vehicles. Finally, the availability of financing options and after-sales service also influence purchasing behavior." "Furthermore, we plan to incorporate user feedback and market trends to continuously improve the tool and ensure it remains at the forefront of presentation technology. Regular updates and feature releases will be essential to maintaining a competitive edge and meeting the evolving needs of our target audience."

# --- Paragraph deletions (content removed by the edit) ---
# Delete from the highest index down so lower indices stay valid as we go.
$p25 = $d.Paragraphs.Item(25)
if ($p25.Range.Text -notmatch [regex]::Escape("The success of these companies will depend on their ability to offer affordable, reliable, and high-performance EVs, coupled with robust charging infrastructure and excellent customer service.")) { Write-Output ("WARNING paragraph 25 text was not as expected before delete: " + $p25.Range.Text) }
$p25.Range.Delete()
$p24 = $d.Paragraphs.Item(24)
if ($p24.Range.Text -notmatch [regex]::Escape("Several new entrants are also emerging, including startups and established automotive companies. These players are focusing on innovative technologies, niche markets, and strategic partnerships. The competitive landscape is expected to intensify in the coming years, with companies vying for market share and customer loyalty.")) { Write-Output ("WARNING paragraph 24 text was not as expected before delete: " + $p24.Range.Text) }
$p24.Range.Delete()
$p23 = $d.Paragraphs.Item(23)
if ($p23.Range.Text -notmatch [regex]::Escape("The Indian EV market is characterized by a mix of domestic and international players. Key competitors include Tata Motors, Mahindra Electric, MG Motor India, and Hyundai Motor India. These companies are offering a range of EV models across different vehicle segments.")) { Write-Output ("WARNING paragraph 23 text was not as expected before delete: " + $p23.Range.Text) }
$p23.Range.Delete()
$p22 = $d.Paragraphs.Item(22)
if ($p22.Range.Text -notmatch [regex]::Escape("Page 2 – Section 3")) { Write-Output ("WARNING paragraph 22 text was not as expected before delete: " + $p22.Range.Text) }
$p22.Range.Delete()
$p17 = $d.Paragraphs.Item(17)
if ($p17.Range.Text -notmatch [regex]::Escape("Moreover, cloud-based platforms for managing charging infrastructure and providing real-time information to EV users are becoming increasingly important. These platforms can optimize charging schedules, monitor battery health, and provide location-based services.")) { Write-Output ("WARNING paragraph 17 text was not as expected before delete: " + $p17.Range.Text) }
$p17.Range.Delete()

# --- Verification summary ---
Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
